$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.224.72'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '3.611.02'
$ws.Range("E3").Value = '  -2.47%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = "'2.32"
$ws.Range("E5").Value = '  +20.88%  '
$ws.Range("D6").Value = "'226.28"
$ws.Range("E6").Value = '  -4.76%  '
$ws.Range("D7").Value = "'636.43"
$ws.Range("E7").Value = '  -2.98%  '
$ws.Range("D8").Value = "'0.411"
$ws.Range("E8").Value = '  -3.74%  '
$ws.Range("D9").Value = "'1.09"
$ws.Range("E9").Value = '  +1.91%  '
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").Value = '3.608.95'
$ws.Range("E11").Value = '  -2.46%  '
$ws.Range("D12").Value = "'46.27"
$ws.Range("E12").Value = '  +4.56%  '
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").Value = "'0.0000290"
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").Value = "'6.46"
$ws.Range("E15").Value = '  -4.49%  '
$ws.Range("D16").Value = '4.285.68'
$ws.Range("E16").Value = '  -2.45%  '
$ws.Range("D17").Value = '95.027.32'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").Value = "'8.75"
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").Value = "'20.20"
$ws.Range("E19").Value = '  +8.34%  '
$ws.Range("D20").Value = '3.598.46'
$ws.Range("E20").Value = '  -2.85%  '
$ws.Range("D21").Value = "'12.84"
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("D22").Value = "'0.514"
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").Value = "'509.48"
$ws.Range("E23").Value = '  -2.69%  '
$ws.Range("D24").Value = "'3.24"
$ws.Range("E24").Value = '  -5.36%  '
$ws.Range("E25").Value = '  +27.57%  '
$ws.Range("D26").Value = "'119.23"
$ws.Range("E26").Value = '  +16.96%  '
$ws.Range("D27").Value = "'0.0000202"
$ws.Range("E27").Value = '  -4.38%  '
$ws.Range("D28").Value = "'6.73"
$ws.Range("E28").Value = '  -2.70%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = "'12.59"
$ws.Range("E29").Value = '  -5.92%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'12.58"
$ws.Range("E30").Value = '  +2.06%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'2.90"
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = "'1.78"
$ws.Range("E34").Value = '  -4.38%  '
$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").Value = "'0.179"
$ws.Range("E35").Value = '  -5.21%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = "'31.69"
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = "'0.583"
$ws.Range("E37").Value = '  -2.64%  '
$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = "'593.89"
$ws.Range("E39").Value = '  -8.18%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = "'8.30"
$ws.Range("E40").Value = '  -6.02%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = "'6.76"
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = "'40.75"
$ws.Range("E42").Value = '  +0.81%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = "'0.158"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.478"
$ws.Range("E44").Value = '  +6.41%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = "'0.0471"
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = "'1.90"
$ws.Range("E46").Value = '  -7.25%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = "'0.915"
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Value = "'23.45"
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'8.53"
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = "'2.20"
$ws.Range("E50").Value = '  -4.01%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = "'53.96"
$ws.Range("E51").Value = '  +0.22%  '
